$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=1;  Col=1; New="48×40="},
    @{Row=1;  Col=2; New="28×87="},
    @{Row=1;  Col=3; New="65×78="},
    @{Row=1;  Col=4; New="46×32="},
    @{Row=1;  Col=5; New="73×91="},
    @{Row=5;  Col=1; New="44×91="},
    @{Row=5;  Col=2; New="35×91="},
    @{Row=5;  Col=3; New="54×47="},
    @{Row=5;  Col=4; New="68×20="},
    @{Row=5;  Col=5; New="87×80="},
    @{Row=10; Col=1; New="50×24="},
    @{Row=10; Col=2; New="44×11="},
    @{Row=10; Col=3; New="46×46="},
    @{Row=10; Col=4; New="35×40="},
    @{Row=10; Col=5; New="90×27="},
    @{Row=15; Col=1; New="87×34="},
    @{Row=15; Col=2; New="86×66="},
    @{Row=15; Col=3; New="83×24="},
    @{Row=15; Col=4; New="13×72="},
    @{Row=15; Col=5; New="61×21="},
    @{Row=20; Col=1; New="89×17="},
    @{Row=20; Col=2; New="13×55="},
    @{Row=20; Col=3; New="85×32="},
    @{Row=20; Col=4; New="44×41="},
    @{Row=20; Col=5; New="71×21="}
)

foreach ($c in $changes) {
    $cell = $t.Cell($c.Row, $c.Col)
    $r = $cell.Range
    $r.MoveEnd(1, -1)
    $r.Text = $c.New
}
